$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.138.96'
$ws.Range('D2').Style = "Normal"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.779.59'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -2.06%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '336.30'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.23%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3894'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.54%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3411'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.64%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.78'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.55%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.182'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -4.32%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07396'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.92%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.003'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('E13').Value = '  -3.91%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.409'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.82%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '1.778.33'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.13%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '7.069'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('E17').Value = '  -3.26%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06657'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '83.23'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -3.82%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '17.49'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.97%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.461'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -1.40%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '27.137.77'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.29%  '
$ws.Range('E24').Value = '  -6.76%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.364'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -4.17%  '
$ws.Range('E26').Value = '  -5.17%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.485'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -7.58%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.434'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -2.06%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '155.44'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +0.87%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.979.49'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -2.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '133.83'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.58%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.973'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -2.57%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.948'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -6.61%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08681'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '12.85'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -7.57%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.614'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -4.47%  '
$ws.Range('E37').Value = '  -4.77%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02367'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.45%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.6746'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.28%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.06320'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.04%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.2180'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -3.48%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.234'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -4.81%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '8.397'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -6.52%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '14.18'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -4.76%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.6347'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.26%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.843'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -3.05%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.120'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.25%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '130.88'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.18%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07100'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.14%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '78.42'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.97%  '
